# Change the MGQ_dict prompt/instructions text for the Metal quiz
# from "Death Metal" wording to "Black Metal" wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MGQ_dict")

# INSTRUCTIONS_METAL (row 3)
$ws.Range("B3").Value = "Sie werden eine Liste mit Name sehen und sollen dort ankreuzen, welcher der Bands eine **BlackMetal Band** ist. Sie haben dazu {{time_out}} Sekunden Zeit."
$ws.Range("C3").Value = "You will be presented with a list of names and you are asked to select all names which are **BlackMetal bands**. You have {{time_out}} seconds to do this."

# PROMPT_METAL (row 4)
$ws.Range("B4").Value = "Bitte wählen Sie alle **Black Metal Bands** aus der untenstehenden Liste aus.  Sie haben {{time_out}} Sekunden Zeit."
$ws.Range("C4").Value = "Please select all  **Black Metal bands**. You have {{time_out}} seconds."

# Match the resulting active-cell selection seen in the saved workbook
$ws.Range("C5").Select()
